# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings that can look like numbers (e.g. "1.004").
# Pre-mark them as Text so Excel stores the literal string instead of silently
# parsing/rounding them into a numeric value (which would drop trailing zeros).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Write the refreshed values.
$ws.Range("D2").Value = "27.042.06"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "1.713.46"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "308.86"
$ws.Range("E5").Value = "  -5.80%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "0.4639"
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("D8").Value = "0.3414"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("D9").Value = "41.84"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "0.07236"
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("D11").Value = "1.038"
$ws.Range("E11").Value = "  -5.42%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "19.74"
$ws.Range("E13").Value = "  -5.50%  "
$ws.Range("D14").Value = "5.813"
$ws.Range("E14").Value = "  -3.63%  "
$ws.Range("D15").Value = "1.719.66"
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("D16").Value = "6.858"
$ws.Range("E16").Value = "  -5.06%  "
$ws.Range("D17").Value = "88.42"
$ws.Range("E17").Value = "  -5.13%  "
$ws.Range("D18").Value = "0.00001033"
$ws.Range("E18").Value = "  -2.47%  "
$ws.Range("D19").Value = "0.06332"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "16.46"
$ws.Range("E21").Value = "  -4.33%  "
$ws.Range("D22").Value = "5.617"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "27.113.20"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").Value = "10.82"
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").Value = "155.29"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").Value = "19.28"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("D28").Value = "1.917.15"
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("D29").Value = "2.108"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").Value = "119.73"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("E31").Value = "  -7.08%  "
$ws.Range("D32").Value = "0.09119"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "3.598"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").Value = "5.315"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("D35").Value = "0.02183"
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("D36").Value = "0.05834"
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("D37").Value = "11.00"
$ws.Range("E37").Value = "  -7.54%  "
$ws.Range("D38").Value = "0.1988"
$ws.Range("E38").Value = "  -5.27%  "
$ws.Range("D39").Value = "4.696"
$ws.Range("E39").Value = "  -5.44%  "
$ws.Range("D40").Value = "1.393"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "0.5889"
$ws.Range("E41").Value = "  -6.60%  "
$ws.Range("D42").Value = "1.126"
$ws.Range("E42").Value = "  -4.70%  "
$ws.Range("D43").Value = "7.431"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").Value = "12.70"
$ws.Range("E44").Value = "  -4.43%  "
$ws.Range("D45").Value = "3.566"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").Value = "0.5613"
$ws.Range("E46").Value = "  -4.34%  "
$ws.Range("D47").Value = "118.95"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").Value = "1.850"
$ws.Range("E48").Value = "  -5.30%  "
$ws.Range("D49").Value = "0.06646"
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("D50").Value = "1.077"
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("E51").Value = "  +0.18%  "
